# Updated symbol list (GitHub Actions data refresh): refresh Price (D) and
# Volume(1h) (E) figures, plus a few re-ranked coin rows (B/C/D/E for rows
# 15-20). Price/Volume cells are stored as plain text in the sheet (not
# numbers/percentages), so values are entered with a leading apostrophe to
# force text, then the style is reset to "Normal" so no stray
# quote-prefix/number-format style is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'29.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'10.13%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.175"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.46%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05707"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.72%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.604"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'2.07%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8590"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'4.36%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.8781"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'3.71%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1367"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'3.21%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.07102"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.73%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.02863"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.78%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09389"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.02%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.001524"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.50%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.04151"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.73%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0006029"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.84%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006156"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.59%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "UpBots"
$ws.Range("C17").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D17").Value = "'0.007491"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'5,108.65%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.480"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.92%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").Value = "'3.039"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.26%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "BTSEToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D20").Value = "'2.186"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-1.81%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.3145"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.01%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.03254"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'4.08%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.1300"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.71%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'3.486"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-2.00%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.1380"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.47%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.005090"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'14.17%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.001217"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-0.03%"
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'23.46%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.03752"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.05%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.005669"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-6.13%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1073"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.85%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002520"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'9.56%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.009797"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'17.31%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005098"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-4.15%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-0.01%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.07099"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-29.70%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002710"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'4.68%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-0.01%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.01%"
$ws.Range("E50").Style = "Normal"
